$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2474043.83
$ws.Range("C7").Value = -44.31689555798313
$ws.Range("D7").Value = 2510
$ws.Range("E7").Value = 2510
$ws.Range("F7").Value = 985.6748326693228
$ws.Range("G7").Value = 5.065809815694
